$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo in B5: remove the stray apostrophe after "akg"
$ws.Range("B5").Value = "축구공 무게는 akg이고 야구공 무게는 bkg이다. 축구공 무게는 야구공 무게의 몇 배인지 구해 보세요."

# Update the active view/selection to reflect the cell edited (B6), with no frozen left column
$ws.Activate()
$ws.Range("B6").Select()
